$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6013534069061279
$ws.Range("B1").Value = 1.434752464294434
$ws.Range("C1").Value = 3.586222410202026
$ws.Range("D1").Value = 0.6105040311813354
$ws.Range("E1").Value = 0.6669391393661499
